# Applies the "adding new mapping from Simon and receiving two new versions
# from Josh" edit to Sheet1 of the workbook:
#   - Row 1/2 header text is refreshed (B2 re-entered).
#   - B3 (old numeric AssayContext_id value) is removed.
#   - The old single data block (rows 3-4) is replaced by three repeated
#     data blocks (rows 3-5, 6-8, 9-11), one per new ad_id value, each
#     carrying the new "purified Salmon sperm DNA" mapping plus the two
#     extra species/tissue rows Josh supplied.
#   - Column widths / selection are refreshed to match the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a plain (non-inferred) text value into a cell without
# introducing any new cell style. We build the text as a formula literal
# and then collapse it down to a plain cached value via copy/paste-values;
# this keeps leading/trailing spaces intact (e.g. "  504727") without
# Excel re-parsing the text as a number. ---
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# ---- Row 1 / Row 2 headers (content unchanged, just re-affirmed) ----
$ws.Cells.Item(1,1).Value = "Source Context"
$ws.Cells.Item(1,5).Value = "Modified Context"

$ws.Cells.Item(2,1).Value = "ad_id"
$ws.Cells.Item(2,2).Value = "AssayContext_id"
$ws.Cells.Item(2,3).Value = "attribue_id"
$ws.Cells.Item(2,4).Value = "value_id"
$ws.Cells.Item(2,5).Value = "attribute_id"
$ws.Cells.Item(2,6).Value = "value_id"
$ws.Cells.Item(2,7).Value = "New_AssayContext_group"

# ---- Clear the old data rows 3 and 4 completely before rebuilding ----
$ws.Range("A3:G4").ClearContents()

# ---- New repeated data blocks ----
$ids = "  504727", "  588364", "  463198"

$row = 3
foreach ($id in $ids) {
    Set-TextValue $ws.Cells.Item($row, 1) $id

    $ws.Cells.Item($row, 3).Value = "biological role"
    $ws.Cells.Item($row, 4).Value = "purified Salmon sperm DNA"
    $ws.Cells.Item($row, 5).Value = "assay component type"
    $ws.Cells.Item($row, 6).Value = "purified DNA"
    $ws.Cells.Item($row, 7).Value = $False

    $ws.Cells.Item($row + 1, 5).Value = "species name"
    $ws.Cells.Item($row + 1, 6).Value = "Salmo salar"
    $ws.Cells.Item($row + 1, 7).Value = $False

    $ws.Cells.Item($row + 2, 5).Value = "biological tissue"
    $ws.Cells.Item($row + 2, 6).Value = "sperm"
    $ws.Cells.Item($row + 2, 7).Value = $False

    $row = $row + 3
}

# ---- Column widths (approximate the refreshed best-fit widths) ----
$ws.Columns.Item(1).ColumnWidth = 6.916666666666666
$ws.Columns.Item(3).ColumnWidth = 12.75
$ws.Columns.Item(4).ColumnWidth = 25.25
$ws.Columns.Item(5).ColumnWidth = 20.083333333333332
$ws.Columns.Item(6).ColumnWidth = 15.750000000000002

# ---- Selection matches the final edited range ----
$ws.Range("C9:G11").Select()
